$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove obsolete columns -------------------------------------------------
# Original layout (A..S):
#  A Norm, Typ                              <- remove
#  B Varumärke                              <- keep -> A
#  C Artikelbenämning                       <- keep -> B
#  D GVM                                    <- keep -> C
#  E Artikelnummer                          <- keep -> D
#  F Typbeteckning                          <- keep -> E
#  G Ritningsnummer                         <- remove
#  H Position                               <- remove
#  I Beteckning                             <- remove
#  J Kompletterande Information övrigt      <- remove
#  K Ref annan                              <- remove
#  L Historiskt Varumärke                   <- remove
#  M Historiskt inköpsreferens              <- remove
#  N Enhet                                  <- keep -> F
#  O Förpackning                            <- remove
#  P SSG-notering                           <- keep -> G
#  Q (empty / unused)                       <- remove
#  R E-nummer                               <- keep -> H
#  S RSK-nummer                             <- keep -> I
# Delete from right to left so earlier column letters stay valid.
$ws.Columns("Q").Delete()
$ws.Columns("O").Delete()
$ws.Columns("G:M").Delete()
$ws.Columns("A").Delete()

# --- Drop the autofilter / sort state ---------------------------------------
$ws.AutoFilterMode = $false

# --- Give the new "SSG-notering" column (now G) a custom width -------------
$ws.Columns("G").ColumnWidth = 28.66666667

# --- Update the selection / active cell to the new last column (H) ---------
$ws.Range("H1:H1048576").Select() | Out-Null

# --- Keep the _FilterDatabase defined name in sync with the new extents ----
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$I`$21217"
    }
}
